$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 13 ---
$ws.Range("A12").Copy()
$ws.Range("A13").PasteSpecial(-4122)
$ws.Range("A13").Value = 41677
$ws.Range("B13").Value = "alterado nome de campo na tabela cashflowsaldo de empresaId para userMaster varchar 255"
$ws.Range("C13").Value = "não"

# --- Row 14 ---
$ws.Range("A12").Copy()
$ws.Range("A14").PasteSpecial(-4122)
$ws.Range("A14").Value = 41677
$ws.Range("B14").Value = "Adicionado campo userMaster nas tabelas cashflowincome e cashflowexpenses"
$ws.Range("C14").Value = "não"

# --- Row 15 ---
$ws.Range("A12").Copy()
$ws.Range("A15").PasteSpecial(-4122)
$ws.Range("A15").Value = 41677
$ws.Range("B15").Value = "Alterado campo userId para userMaster varchar 255"
$ws.Range("C15").Value = "não"

$excel.CutCopyMode = $false

$ws.Range("B17").Select()
